# Apply updated crypto price / volume(1h) figures (rows 2-51, columns D & E).
# Values are written as text (apostrophe-prefixed) to match the original
# inlineStr cell type, then Style is reset to "Normal" so no numeric
# formatting / style index is introduced by Excel's text auto-detection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'88.081.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.85%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.113.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.22%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'215.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.25%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'635.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.86%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.390"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.90%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +17.67%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.09%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'3.110.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.24%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.566"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.29%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.97%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.39%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.16%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'87.869.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.84%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'3.683.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.77%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'32.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.39%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.105.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.41%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'3.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.25%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  +20.61%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'13.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'423.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.03%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'8.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.34%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'4.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.50%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'5.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +8.45%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'82.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +11.23%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'11.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.51%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'3.283.70"
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -0.04%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.155"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -6.19%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'4.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.88%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'8.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.75%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.150"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +19.04%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'502.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.28%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'6.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.22%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.05%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -0.08%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'22.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.07%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'22.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.24%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.44%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D43").Value = "'0.365"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.40%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  -2.25%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.135"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +10.23%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'146.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.06%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'43.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.52%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'162.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.11%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0647"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +11.46%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.719"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.63%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -2.09%  "
$ws.Range("E51").Style = "Normal"
